# Ads1299_defRegs.xlsx edit
# - CONFIG2 (row 3): G3 restyled to the "Neutral" centered look used
#   elsewhere in the bit columns (copy format from G5, value untouched);
#   J3 bit flips 1->0, K3 bit flips 0->1. C3's shared formula recalculates
#   0xD2 -> 0xD1.
# - CONFIG3 (row 5): J5 bit flips 1->0. C5 recalculates 0x02 -> 0x00.
# - row 8: I8 bit flips 0->1. C8 recalculates 0xE1 -> 0xE5.
# - selection moves from D3 down to N31 (scrolled into the notes column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3 picks up the same visual style (fill/font/centered alignment) already
# used on the sibling bit cell G5, without touching G3's value.
$ws.Range("G5").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# Flip the individual register bits.
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("J5").Value = 0
$ws.Range("I8").Value = 1

# Move the active selection to N31, matching the committed view state.
$ws.Range("N31").Select() | Out-Null
